# Applies the cell-value changes from the commit diff (crypto price/volume refresh).
# Values are written with a leading apostrophe + Style reset so that numeric-looking
# strings (e.g. "214.43", "0.0790", "1.00") are stored as literal text, matching the
# original workbook (every data cell is an inline/shared text string, never a Number).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.867.22"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "1.632.37"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  +0.65%  "
$ws.Range("D5").Value = "'214.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").Value = "'0.0633"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").Value = "'19.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.55%  "
$ws.Range("D11").Value = "'0.0790"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").Value = "1.664.51"
$ws.Range("E12").Value = "  +2.47%  "
$ws.Range("D13").Value = "1.857.29"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").Value = "'62.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").Value = "25.879.58"
$ws.Range("E18").Value = "  +0.22%  "
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'4.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'193.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("D22").Value = "'9.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D24").Value = "'1.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("D25").Value = "'143.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.61%  "
$ws.Range("E27").Value = "  +3.18%  "
$ws.Range("D28").Value = "'6.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("D30").Value = "'1.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("D31").Value = "'0.0498"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.85%  "
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("D33").Value = "'3.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("E34").Value = "  -1.75%  "
$ws.Range("E35").Value = "  +1.52%  "
$ws.Range("D36").Value = "'0.900"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("D37").Value = "1.136.33"
$ws.Range("E37").Value = "  -0.48%  "
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("E40").Value = "  +0.25%  "
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'5.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.37%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'99.05"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.52%  "
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("D45").Value = "1.766.77"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("E46").Value = "  +2.93%  "
$ws.Range("D47").Value = "'56.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.61%  "
$ws.Range("D48").Value = "'0.0526"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.85%  "
$ws.Range("E49").Value = "  -0.96%  "
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("E51").Value = "  +0.59%  "
